# Add a new worksheet "ODI Bowling Extra" (mirrors the existing "ODI Batting
# Extra" sheet, but for bowling: maiden overs + % of all wickets taken by
# match) as the 5th/last sheet in the workbook.

$wb = $excel.ActiveWorkbook

# Copy the header cell formatting (bold, centered, thin border) from the
# existing "ODI Batting Extra" sheet so the new sheet's header matches the
# look of the other "Extra" sheet.
$sourceSheet = $wb.Worksheets.Item("ODI Batting Extra")
$sourceSheet.Range("A1").Copy()

# Insert the new sheet after the last existing sheet ("ODI Batting Extra")
# so tab order stays Player Info, ODI Batting, ODI Bowling, ODI Batting
# Extra, ODI Bowling Extra.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Sheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Bowling Extra"

# Header row
$headerRange = $newSheet.Range("A1:C1")
$headerRange.PasteSpecial(-4122) # xlPasteFormats

$newSheet.Cells.Item(1, 1).Value = "MATCH_CODE"
$newSheet.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$newSheet.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"

# Force the data area to be stored as text (matches the source data, which
# keeps numeric-looking values like maiden-over counts and percentages as
# plain text rather than numbers).
$dataRange = $newSheet.Range("A2:C21")
$dataRange.NumberFormat = "@"

# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$data = @(
    @("3519", "", ""),
    @("3602", "0", ""),
    @("3606", "", ""),
    @("3608", "2", "30.00%"),
    @("3622", "1", ""),
    @("3625", "0", "10.00%"),
    @("3629", "0", "20.00%"),
    @("3727", "0", ""),
    @("3735", "0", ""),
    @("3738", "", ""),
    @("3744", "1", "20.00%"),
    @("3746", "", ""),
    @("3749", "0", "20.00%"),
    @("3756", "0", ""),
    @("3761", "", ""),
    @("3769", "1", ""),
    @("3780", "0", "10.00%"),
    @("3785", "1", "10.00%"),
    @("3890", "0", "10.00%"),
    @("3891", "1", "")
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}


# Restore the originally-active sheet (adding a sheet makes it active by
# default in Excel) so the workbook's selected-tab state is left as found.
$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select()
